$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few text cells hold values that LOOK numeric/date-like ("9.27", "9/1/7143", ...).
# Left alone, Excel would auto-convert these to a Double/date serial on assignment,
# so force the cell to Text format first, then restore the default "Normal" style
# afterwards so no stray number-format style is left behind on the cell.
$forceTextAddrs = @("I4", "D10", "I16", "I17")
foreach ($addr in $forceTextAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 4
$ws.Range("B4").Value = "4-ANPP"
$ws.Range("C4").Value = 88890
$ws.Range("D4").Value = "21409-26-7"
$ws.Range("E4").Value = "Q88EHD0U8G"
$ws.Range("F4").Value = ""
$ws.Range("I4").Value = "9.27"
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0
$ws.Range("AJ4").Value = 1
$ws.Range("AK4").Value = 1
# Row 5
$ws.Range("B5").Value = "ethyl-4-ANPP"
$ws.Range("C5").Value = 156346345
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "trace"
$ws.Range("I5").Value = "."
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1
$ws.Range("AJ5").Value = 0
$ws.Range("AK5").Value = 0
# Row 9
$ws.Range("B9").Value = "ketamine"
$ws.Range("C9").Value = 3821
$ws.Range("D9").Value = "6740-88-1"
$ws.Range("E9").Value = "690G0D6V8H"
$ws.Range("AF9").Value = 1
# Row 10
$ws.Range("B10").Value = "ecgonine methylester (EME)"
$ws.Range("C10").Value = 104904
$ws.Range("D10").Value = "9/1/7143"
$ws.Range("E10").Value = "Y35FJB3QBJ"
$ws.Range("AF10").Value = 0
# Row 14
$ws.Range("B14").Value = "xylazine"
$ws.Range("C14").Value = 5707
$ws.Range("D14").Value = "7361-61-7"
$ws.Range("E14").Value = "2KFG9TP5V8"
$ws.Range("F14").Value = "trace"
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("S14").Value = 1
$ws.Range("T14").Value = 1
$ws.Range("AJ14").Value = 0
# Row 15
$ws.Range("B15").Value = "fentanyl"
$ws.Range("C15").Value = 3345
$ws.Range("D15").Value = "437-38-7"
$ws.Range("E15").Value = "UF599785JZ"
$ws.Range("F15").Value = ""
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 0
$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 1
$ws.Range("S15").Value = 0
$ws.Range("T15").Value = 0
$ws.Range("AJ15").Value = 1
# Row 16
$ws.Range("B16").Value = "dimethyl sulfone (methylsulfonylmethane MSM)"
$ws.Range("C16").Value = 6213
$ws.Range("D16").Value = "67-71-0"
$ws.Range("E16").Value = "9H4PO4Z4FT"
$ws.Range("F16").Value = ""
$ws.Range("I16").Value = "2.38"
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 0
# Row 17
$ws.Range("B17").Value = "ketamine"
$ws.Range("C17").Value = 3821
$ws.Range("D17").Value = "6740-88-1"
$ws.Range("E17").Value = "690G0D6V8H"
$ws.Range("I17").Value = "7.21"
$ws.Range("AF17").Value = 1
$ws.Range("AG17").Value = 1
# Row 18
$ws.Range("B18").Value = "N,N-dimethyltryptamine (DMT)"
$ws.Range("C18").Value = 6089
$ws.Range("D18").Value = "61-50-7"
$ws.Range("E18").Value = "WUB601BHAA"
$ws.Range("F18").Value = "trace"
$ws.Range("I18").Value = "."
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1
$ws.Range("AF18").Value = 0
$ws.Range("AG18").Value = 0

# Restore default styling on the forced-text cells.
foreach ($addr in $forceTextAddrs) {
    $ws.Range($addr).Style = "Normal"
}
